# Update the "cryptos" price/volume snapshot on Sheet1.
# Numeric-looking price strings (e.g. "214.44") are prefixed with a leading
# apostrophe so Excel stores them as text, matching the original inlineStr
# cell type instead of being auto-converted to a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.910.77'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '1.626.55'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('D5').Value = '''214.44'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').Value = '''29.73'
$ws.Range('E8').Value = '  +9.26%  '
$ws.Range('E9').Value = '  +2.56%  '
$ws.Range('D10').Value = '''0.0612'
$ws.Range('E10').Value = '  +1.71%  '
$ws.Range('D11').Value = '''0.0914'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '1.858.44'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').Value = '1.627.01'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').Value = '''0.573'
$ws.Range('E14').Value = '  +6.43%  '
$ws.Range('E15').Value = '  +4.75%  '
$ws.Range('D16').Value = '29.971.93'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '''8.91'
$ws.Range('E17').Value = '  +17.18%  '
$ws.Range('D18').Value = '''64.75'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').Value = '''244.11'
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('E22').Value = '  +3.36%  '
$ws.Range('D23').Value = '''9.64'
$ws.Range('E23').Value = '  +4.43%  '
$ws.Range('D24').Value = '''2.13'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('D25').Value = '''157.61'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('E28').Value = '  +2.86%  '
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('D30').Value = '''0.0490'
$ws.Range('E30').Value = '  +3.10%  '
$ws.Range('E31').Value = '  +5.52%  '
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').Value = '1.422.79'
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('D35').Value = '''1.63'
$ws.Range('E35').Value = '  +5.99%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '''2.86'
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('E39').Value = '  +2.99%  '
$ws.Range('D40').Value = '''0.560'
$ws.Range('E40').Value = '  +3.44%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '''1.99'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''0.834'
$ws.Range('E42').Value = '  +4.04%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '''0.0501'
$ws.Range('E43').Value = '  +2.15%  '
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = '''69.49'
$ws.Range('E45').Value = '  +5.18%  '
$ws.Range('E46').Value = '  +11.28%  '
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').Value = '''5.40'
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('D49').Value = '1.766.71'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').Value = '''88.95'
$ws.Range('E50').Value = '  +2.37%  '
$ws.Range('D51').Value = '0.0₆0109'
$ws.Range('E51').Value = '  +2.96%  '
